$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I2").Value = 0.07819089147286823
$ws.Range("J2").Value = 0.07819089147286823
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 0.1020818183178889
$ws.Range("R2").Value = 0.9187363648610001
$ws.Range("S2").Value = 0.0008164800444740656
$ws.Range("T2").Value = 0.0008164800444740656
$ws.Range("I3").Value = 0.07819089147286823
$ws.Range("J3").Value = 0.07819089147286823
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 6.546257245869779
$ws.Range("R3").Value = 58.91631521282801
$ws.Range("S3").Value = 0.05235886757622329
$ws.Range("T3").Value = 0.05235886757622329
$ws.Range("I4").Value = 0.07819089147286823
$ws.Range("J4").Value = 0.07819089147286823
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 0.8015873433676667
$ws.Range("R4").Value = 7.214286090309
$ws.Range("S4").Value = 0.006411328486769214
$ws.Range("T4").Value = 0.006411328486769213
$ws.Range("I5").Value = 0.07819089147286823
$ws.Range("J5").Value = 0.07819089147286823
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 2.326023943550445
$ws.Range("R5").Value = 20.934215491954
$ws.Range("S5").Value = 0.01860421536540165
$ws.Range("T5").Value = 0.01860421536540165
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3171023333333333
$ws.Range("H6").Value = 0.9513069999999999
$ws.Range("I6").Value = 0.9218091085271318
$ws.Range("J6").Value = 0.9218091085271318
$ws.Range("M6").Value = 3.795192333333334
$ws.Range("N6").Value = 11.385577
$ws.Range("O6").Value = 0.01044213755712683
$ws.Range("P6").Value = 0.01044213755712683
$ws.Range("Q6").Value = 1.203464344348778
$ws.Range("R6").Value = 10.831179099139
$ws.Range("S6").Value = 0.009625657512652768
$ws.Range("T6").Value = 0.009625657512652768
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3171023333333333
$ws.Range("H7").Value = 0.9513069999999999
$ws.Range("I7").Value = 0.9218091085271318
$ws.Range("J7").Value = 0.9218091085271318
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.6696287328350964
$ws.Range("P7").Value = 0.6696287328350964
$ws.Range("Q7").Value = 77.17522389546356
$ws.Range("R7").Value = 694.577015059172
$ws.Range("S7").Value = 0.6172698652588732
$ws.Range("T7").Value = 0.6172698652588732
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3171023333333333
$ws.Range("H8").Value = 0.9513069999999999
$ws.Range("I8").Value = 0.9218091085271318
$ws.Range("J8").Value = 0.9218091085271318
$ws.Range("M8").Value = 29.801371
$ws.Range("N8").Value = 89.404113
$ws.Range("O8").Value = 0.08199584844219236
$ws.Range("P8").Value = 0.08199584844219235
$ws.Range("Q8").Value = 9.450084280632332
$ws.Range("R8").Value = 85.05075852569098
$ws.Range("S8").Value = 0.07558451995542315
$ws.Range("T8").Value = 0.07558451995542315
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3171023333333333
$ws.Range("H9").Value = 0.9513069999999999
$ws.Range("I9").Value = 0.9218091085271318
$ws.Range("J9").Value = 0.9218091085271318
$ws.Range("M9").Value = 86.47679266666667
$ws.Range("N9").Value = 259.430378
$ws.Range("O9").Value = 0.2379332811655844
$ws.Range("P9").Value = 0.2379332811655844
$ws.Range("Q9").Value = 27.42199273378289
$ws.Range("R9").Value = 246.797934604046
$ws.Range("S9").Value = 0.2193290658001827
$ws.Range("T9").Value = 0.2193290658001827
